$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value would otherwise be auto-converted to a number
# by type inference; force them to stay text (the source workbook stores
# every Price/Volume cell as an inline string).
$textCells = @('D5', 'D6', 'D8', 'D14', 'D17', 'D20', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D30', 'D31', 'D32', 'D33', 'D37', 'D38', 'D39', 'D40', 'D41', 'D43', 'D44', 'D45', 'D47', 'D48', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.502.41'
$ws.Range('E2').Value = '  -1.29%  '

$ws.Range('D3').Value = '3.022.72'
$ws.Range('E3').Value = '  -1.22%  '

$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').Value = '586.45'
$ws.Range('E5').Value = '  -0.68%  '

$ws.Range('D6').Value = '148.07'
$ws.Range('E6').Value = '  -3.28%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '0.526'
$ws.Range('E8').Value = '  -2.16%  '

$ws.Range('D9').Value = '3.024.35'
$ws.Range('E9').Value = '  -1.20%  '

$ws.Range('E10').Value = '  -3.77%  '

$ws.Range('E11').Value = '  -1.03%  '

$ws.Range('E12').Value = '  -1.37%  '

$ws.Range('E13').Value = '  -2.73%  '

$ws.Range('D14').Value = '34.71'
$ws.Range('E14').Value = '  -4.82%  '

$ws.Range('E15').Value = '  +2.08%  '

$ws.Range('D16').Value = '3.521.75'
$ws.Range('E16').Value = '  -1.20%  '

$ws.Range('D17').Value = '7.08'
$ws.Range('E17').Value = '  -1.75%  '

$ws.Range('D18').Value = '62.387.43'
$ws.Range('E18').Value = '  -1.46%  '

$ws.Range('D19').Value = '3.021.10'
$ws.Range('E19').Value = '  -1.36%  '

$ws.Range('D20').Value = '463.15'
$ws.Range('E20').Value = '  -4.50%  '

$ws.Range('E21').Value = '  -3.63%  '

$ws.Range('D22').Value = '0.686'
$ws.Range('E22').Value = '  -2.75%  '

$ws.Range('D23').Value = '7.46'
$ws.Range('E23').Value = '  -0.58%  '

$ws.Range('D24').Value = '81.61'
$ws.Range('E24').Value = '  -0.59%  '

$ws.Range('D25').Value = '2.27'
$ws.Range('E25').Value = '  -5.38%  '

$ws.Range('D26').Value = '12.43'
$ws.Range('E26').Value = '  -2.58%  '

$ws.Range('D27').Value = '10.20'
$ws.Range('E27').Value = '  -2.44%  '

$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.00%  '

$ws.Range('E29').Value = '  -0.22%  '

$ws.Range('D30').Value = '2.63'
$ws.Range('E30').Value = '  -1.91%  '

$ws.Range('D31').Value = '7.16'
$ws.Range('E31').Value = '  -4.09%  '

$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '2.11'
$ws.Range('E32').Value = '  -5.24%  '

$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '28.47'
$ws.Range('E33').Value = '  +4.32%  '

$ws.Range('E34').Value = '  -1.83%  '

$ws.Range('D35').Value = '0.0₃0808'
$ws.Range('E35').Value = '  -1.53%  '

$ws.Range('E36').Value = '  -2.36%  '

$ws.Range('D37').Value = '5.79'
$ws.Range('E37').Value = '  -3.88%  '

$ws.Range('D38').Value = '2.13'
$ws.Range('E38').Value = '  -3.79%  '

$ws.Range('D39').Value = '50.49'
$ws.Range('E39').Value = '  -0.14%  '

$ws.Range('D40').Value = '9.10'
$ws.Range('E40').Value = '  -2.14%  '

$ws.Range('D41').Value = '2.95'
$ws.Range('E41').Value = '  -8.82%  '

$ws.Range('E42').Value = '  +0.33%  '

$ws.Range('D43').Value = '393.97'
$ws.Range('E43').Value = '  -10.09%  '

$ws.Range('D44').Value = '0.277'
$ws.Range('E44').Value = '  -3.67%  '

$ws.Range('D45').Value = '0.0358'
$ws.Range('E45').Value = '  -1.05%  '

$ws.Range('D46').Value = '2.760.90'
$ws.Range('E46').Value = '  -1.97%  '

$ws.Range('D47').Value = '37.47'
$ws.Range('E47').Value = '  -4.97%  '

$ws.Range('D48').Value = '128.77'
$ws.Range('E48').Value = '  -3.24%  '

$ws.Range('E49').Value = '  +0.08%  '

$ws.Range('D50').Value = '0.109'
$ws.Range('E50').Value = '  -0.69%  '

$ws.Range('D51').Value = '24.23'
$ws.Range('E51').Value = '  -4.07%  '
